$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Bump version + date (row 3 = Version, row 8 = Date)
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" metadata row right after "Contact" (row 10), before "Description" (row 11)
$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Row 5 = Component.typeId -> add the II-1 constraint text in the "Constraint(s)" column (AJ)
$elements.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}" + [char]10
